# Update crypto price/volume table cells to reflect refreshed data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "61.023.41"
$ws.Range("E2").Value = "  +1.02%  "
Set-TextValue $ws.Range("D3") "3.384.54"
$ws.Range("E3").Value = "  +0.10%  "
Set-TextValue $ws.Range("D5") "571.38"
$ws.Range("E5").Value = "  -0.03%  "
Set-TextValue $ws.Range("D6") "141.51"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  -1.64%  "
Set-TextValue $ws.Range("D12") "3.962.44"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  +2.12%  "
Set-TextValue $ws.Range("D14") "27.87"
$ws.Range("E14").Value = "  -1.17%  "
Set-TextValue $ws.Range("D15") "3.381.38"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  -0.04%  "
Set-TextValue $ws.Range("D17") "61.108.48"
$ws.Range("E17").Value = "  +0.92%  "
Set-TextValue $ws.Range("D18") "6.12"
$ws.Range("E18").Value = "  -2.37%  "
Set-TextValue $ws.Range("D19") "13.63"
$ws.Range("E19").Value = "  -3.32%  "
Set-TextValue $ws.Range("D20") "8.94"
$ws.Range("E20").Value = "  -3.04%  "
Set-TextValue $ws.Range("D21") "384.11"
$ws.Range("E21").Value = "  -1.20%  "
Set-TextValue $ws.Range("D22") "75.59"
$ws.Range("E22").Value = "  +2.79%  "
Set-TextValue $ws.Range("D23") "0.552"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("E24").Value = "  -0.12%  "
Set-TextValue $ws.Range("D25") "0.0000116"
$ws.Range("E25").Value = "  -1.04%  "
Set-TextValue $ws.Range("D26") "3.521.68"
$ws.Range("E26").Value = "  +0.10%  "
Set-TextValue $ws.Range("D27") "0.185"
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -2.36%  "
Set-TextValue $ws.Range("D30") "7.98"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  -4.32%  "
Set-TextValue $ws.Range("D34") "23.19"
$ws.Range("E34").Value = "  -2.36%  "
Set-TextValue $ws.Range("D35") "6.95"
$ws.Range("E35").Value = "  +0.20%  "
Set-TextValue $ws.Range("D36") "165.97"
$ws.Range("E36").Value = "  -0.79%  "
Set-TextValue $ws.Range("D37") "3.418.74"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  +0.58%  "
Set-TextValue $ws.Range("D39") "1.47"
$ws.Range("E39").Value = "  -2.41%  "
Set-TextValue $ws.Range("D40") "0.0766"
$ws.Range("E40").Value = "  -1.37%  "
Set-TextValue $ws.Range("D41") "26.61"
$ws.Range("E41").Value = "  -1.65%  "
Set-TextValue $ws.Range("D42") "1.00"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("E46").Value = "  +0.20%  "
Set-TextValue $ws.Range("D47") "2.449.92"
$ws.Range("E47").Value = "  -3.31%  "
Set-TextValue $ws.Range("D48") "22.92"
$ws.Range("E48").Value = "  -0.81%  "
Set-TextValue $ws.Range("D49") "6.66"
$ws.Range("E49").Value = "  -2.74%  "
Set-TextValue $ws.Range("D50") "2.14"
$ws.Range("E50").Value = "  +10.17%  "
Set-TextValue $ws.Range("D51") "0.0262"
$ws.Range("E51").Value = "  -1.86%  "
